# "modify design decision slide"
# The "Design Decisions" slide's content placeholder has a paragraph that
# reads "Received messages are shown in green". Update the leading word
# "Received" to "Received (read)" while leaving the rest of the paragraph
# (and its existing run/formatting split) untouched.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -like "*Received messages are shown in green*") {
                $targetSlide = $s
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$startPos = $fullText.IndexOf("Received messages are shown in green") + 1
$word = $tr.Characters($startPos, 8)
$word.Text = "Received (read)"

Write-Output ("Updated paragraph now reads: " + $targetShape.TextFrame.TextRange.Text)
